$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 80, column A: the new date label "05-12-2025".
# NOTE: assigning this literal date-like string straight to .Value (or
# .Value2/.Formula) makes Excel's type-inference treat it as a real date
# serial and stamp a brand-new number-format style on the cell — which does
# not match the source data (every other date in column A is plain text
# sharing style index 3). To keep the literal text *and* the original
# style, compute it with TRIM() of a space-padded literal (space prefix
# defeats the date auto-detection) and then flatten the formula down to a
# static value with Copy/PasteSpecial (xlPasteValues = -4163), which does
# not touch the cell's number format/style.
$ws.Cells.Item(80, 1).Formula = "=TRIM(`" 05-12-2025`")"
$ws.Cells.Item(80, 1).Copy()
$ws.Cells.Item(80, 1).PasteSpecial(-4163)

# Row 80, column B: the new gold-price sentence. Plain text (doesn't look
# like a number/date/bool) so a direct .Value assignment stores it as a
# literal shared string and keeps the inherited column style (index 4).
$ws.Cells.Item(80, 2).Value = "The price of gold in India today is ₹12,993 per gram for 24 karat gold, ₹11,910 per gram for 22 karat gold and ₹9,745 per gram for 18 karat gold (also called 999 gold)."
